$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old regression test-case rows (rows 3-6), shifting remaining
# rows up so the header/used-range collapses to A1:C2.
$ws.Range("A3:C6").EntireRow.Delete()

# Replace the remaining data row with the newly added iAuthor test case.
$ws.Range("A2").Value = "iAU_TC_ID_131"
$ws.Range("B2").Value = "@RegressionA Validation of Create Exam – Start from scratch(Negative Scenario)"
$ws.Range("C2").Value = "passed"
